# Apply updated crypto price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '60.887.34'
$ws.Range('E2').Value = '  -2.59%  '

# Row 3
$ws.Range('D3').Value = '3.350.50'
$ws.Range('E3').Value = '  -2.46%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.00%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '566.67'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.94%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.95'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.85%  '

# Row 7
$ws.Range('E7').Value = '  +0.09%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.484'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.69%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.92'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.18%  '

# Row 10
$ws.Range('E10').Value = '  -1.02%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.414'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.90%  '

# Row 12
$ws.Range('D12').Value = '3.927.41'
$ws.Range('E12').Value = '  -2.46%  '

# Row 13
$ws.Range('E13').Value = '  +0.82%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.72'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.62%  '

# Row 15
$ws.Range('D15').Value = '3.350.70'
$ws.Range('E15').Value = '  -2.74%  '

# Row 16
$ws.Range('E16').Value = '  -1.11%  '

# Row 17
$ws.Range('D17').Value = '60.860.10'
$ws.Range('E17').Value = '  -2.72%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.31'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.75%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.43'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.63%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.86'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.86%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '375.97'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.57%  '

# Row 22
$ws.Range('E22').Value = '  +0.08%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '74.88'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.27%  '

# Row 24
$ws.Range('E24').Value = '  -0.07%  '

# Row 25
$ws.Range('D25').Value = '3.500.00'
$ws.Range('E25').Value = '  -2.31%  '

# Row 26
$ws.Range('E26').Value = '  -5.24%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.174'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.16%  '

# Row 28
$ws.Range('E28').Value = '  +0.19%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.34'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.09%  '

# Row 30
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.08'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.87%  '

# Row 31
$ws.Range('B31').Value = 'USDe'
$ws.Range('C31').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.01%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.71'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.03%  '

# Row 33
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '22.83'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.35%  '

# Row 34
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.30'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.65%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.32'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.31%  '

# Row 36
$ws.Range('E36').Value = '  -0.88%  '

# Row 37
$ws.Range('E37').Value = '  -4.12%  '

# Row 38
$ws.Range('E38').Value = '  -2.02%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '28.83'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -8.98%  '

# Row 40
$ws.Range('D40').Value = '3.388.06'
$ws.Range('E40').Value = '  -2.40%  '

# Row 41
$ws.Range('E41').Value = '  -2.47%  '

# Row 42
$ws.Range('E42').Value = '  -3.31%  '

# Row 43
$ws.Range('E43').Value = '  -1.00%  '

# Row 44
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.61'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.02%  '

# Row 45
$ws.Range('B45').Value = 'ONDO'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.13'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.31%  '

# Row 46
$ws.Range('D46').Value = '2.468.54'
$ws.Range('E46').Value = '  -3.63%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.65'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.76%  '

# Row 48
$ws.Range('B48').Value = 'FirstDigitalUSD'
$ws.Range('C48').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.999'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.10%  '

# Row 49
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '22.34'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.92%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0260'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.07%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.816'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.76%  '
